$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (max/possible-points row): add new Assignment 7 max value
$ws.Range("J2").Value = 20

# Row 7 - Covell, David A.: add assignment scores (Assignment 1-7)
$ws.Range("D7").Value = 41
$ws.Range("E7").Value = 19
$ws.Range("F7").Value = 24
$ws.Range("G7").Value = 13
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 20
$ws.Range("J7").Value = 20

# Row 9 - Davis, Troy W.: add Assignment 6 score
$ws.Range("I9").Value = 20

# Row 11 - Estrada, Andres A.: add Assignment 1 score
$ws.Range("D11").Value = 33

# Row 18 - Le, Jimmy: add assignment scores (Assignment 1-5)
$ws.Range("D18").Value = 41
$ws.Range("E18").Value = 19
$ws.Range("F18").Value = 24
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 10

# Update selection to match the final cursor position recorded in the workbook
$ws.Range("F11").Select()
